$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 45.01222466666667
$ws.Range("H2").Value = 135.036674
$ws.Range("I2").Value = 0.7482903203664146
$ws.Range("J2").Value = 0.7482903203664146
$ws.Range("K2").Value = 3.0
$ws.Range("L2").Value = 1.0
$ws.Range("M2").Value = 3.825035
$ws.Range("N2").Value = 11.475105
$ws.Range("O2").Value = 0.03111562857396839
$ws.Range("P2").Value = 0.03111562857396839
$ws.Range("Q2").Value = 172.1733347778633
$ws.Range("R2").Value = 1549.56001300077
$ws.Range("S2").Value = 0.02328352367401717
$ws.Range("T2").Value = 0.02328352367401717
$ws.Range("G3").Value = 45.01222466666667
$ws.Range("H3").Value = 135.036674
$ws.Range("I3").Value = 0.7482903203664146
$ws.Range("J3").Value = 0.7482903203664146
$ws.Range("O3").Value = 0.4709815605157605
$ws.Range("P3").Value = 0.4709815605157605
$ws.Range("Q3").Value = 2606.100844150117
$ws.Range("R3").Value = 23454.90759735106
$ws.Range("S3").Value = 0.3524309428050123
$ws.Range("T3").Value = 0.3524309428050123
$ws.Range("G4").Value = 45.01222466666667
$ws.Range("H4").Value = 135.036674
$ws.Range("I4").Value = 0.7482903203664146
$ws.Range("J4").Value = 0.7482903203664146
$ws.Range("M4").Value = 61.10114166666667
$ws.Range("N4").Value = 183.303425
$ws.Range("O4").Value = 0.4970413158429724
$ws.Range("P4").Value = 0.4970413158429724
$ws.Range("Q4").Value = 2750.298316089828
$ws.Range("R4").Value = 24752.68484480845
$ws.Range("S4").Value = 0.3719312054674821
$ws.Range("T4").Value = 0.3719312054674821
$ws.Range("G5").Value = 45.01222466666667
$ws.Range("H5").Value = 135.036674
$ws.Range("I5").Value = 0.7482903203664146
$ws.Range("J5").Value = 0.7482903203664146
$ws.Range("M5").Value = 0.1059033333333333
$ws.Range("N5").Value = 0.31771
$ws.Range("O5").Value = 0.0008614950672987739
$ws.Range("P5").Value = 0.0008614950672987739
$ws.Range("Q5").Value = 4.766944632948889
$ws.Range("R5").Value = 42.90250169654
$ws.Range("S5").Value = 0.0006446484199030854
$ws.Range("T5").Value = 0.0006446484199030854
$ws.Range("G6").Value = 2.766295666666667
$ws.Range("H6").Value = 8.298887
$ws.Range("I6").Value = 0.04598733535094824
$ws.Range("J6").Value = 0.04598733535094825
$ws.Range("K6").Value = 3.0
$ws.Range("L6").Value = 1.0
$ws.Range("M6").Value = 3.825035
$ws.Range("N6").Value = 11.475105
$ws.Range("O6").Value = 0.03111562857396839
$ws.Range("P6").Value = 0.03111562857396839
$ws.Range("Q6").Value = 10.58117774534833
$ws.Range("R6").Value = 95.230599708135
$ws.Range("S6").Value = 0.001430924845886631
$ws.Range("T6").Value = 0.001430924845886632
$ws.Range("G7").Value = 2.766295666666667
$ws.Range("H7").Value = 8.298887
$ws.Range("I7").Value = 0.04598733535094824
$ws.Range("J7").Value = 0.04598733535094825
$ws.Range("O7").Value = 0.4709815605157605
$ws.Range("P7").Value = 0.4709815605157605
$ws.Range("Q7").Value = 160.1619454593975
$ws.Range("R7").Value = 1441.457509134577
$ws.Range("S7").Value = 0.0216591869675512
$ws.Range("T7").Value = 0.0216591869675512
$ws.Range("G8").Value = 2.766295666666667
$ws.Range("H8").Value = 8.298887
$ws.Range("I8").Value = 0.04598733535094824
$ws.Range("J8").Value = 0.04598733535094825
$ws.Range("M8").Value = 61.10114166666667
$ws.Range("N8").Value = 183.303425
$ws.Range("O8").Value = 0.4970413158429724
$ws.Range("P8").Value = 0.4970413158429724
$ws.Range("Q8").Value = 169.0238234208861
$ws.Range("R8").Value = 1521.214410787975
$ws.Range("S8").Value = 0.02285760567494735
$ws.Range("T8").Value = 0.02285760567494736
$ws.Range("G9").Value = 2.766295666666667
$ws.Range("H9").Value = 8.298887
$ws.Range("I9").Value = 0.04598733535094824
$ws.Range("J9").Value = 0.04598733535094825
$ws.Range("M9").Value = 0.1059033333333333
$ws.Range("N9").Value = 0.31771
$ws.Range("O9").Value = 0.0008614950672987739
$ws.Range("P9").Value = 0.0008614950672987739
$ws.Range("Q9").Value = 0.2929599320855555
$ws.Range("R9").Value = 2.63663938877
$ws.Range("S9").Value = 0.00003961786256305644
$ws.Range("T9").Value = 0.00003961786256305644
$ws.Range("G10").Value = 12.37490333333333
$ws.Range("H10").Value = 37.12471
$ws.Range("I10").Value = 0.2057223442826371
$ws.Range("J10").Value = 0.2057223442826371
$ws.Range("K10").Value = 3.0
$ws.Range("L10").Value = 1.0
$ws.Range("M10").Value = 3.825035
$ws.Range("N10").Value = 11.475105
$ws.Range("O10").Value = 0.03111562857396839
$ws.Range("P10").Value = 0.03111562857396839
$ws.Range("Q10").Value = 47.33443837161666
$ws.Range("R10").Value = 426.00994534455
$ws.Range("S10").Value = 0.006401180054064586
$ws.Range("T10").Value = 0.006401180054064586
$ws.Range("G11").Value = 12.37490333333333
$ws.Range("H11").Value = 37.12471
$ws.Range("I11").Value = 0.2057223442826371
$ws.Range("J11").Value = 0.2057223442826371
$ws.Range("O11").Value = 0.4709815605157605
$ws.Range("P11").Value = 0.4709815605157605
$ws.Range("Q11").Value = 716.4774961047123
$ws.Range("R11").Value = 6448.29746494241
$ws.Range("S11").Value = 0.09689143074319698
$ws.Range("T11").Value = 0.09689143074319698
$ws.Range("G12").Value = 12.37490333333333
$ws.Range("H12").Value = 37.12471
$ws.Range("I12").Value = 0.2057223442826371
$ws.Range("J12").Value = 0.2057223442826371
$ws.Range("M12").Value = 61.10114166666667
$ws.Range("N12").Value = 183.303425
$ws.Range("O12").Value = 0.4970413158429724
$ws.Range("P12").Value = 0.4970413158429724
$ws.Range("Q12").Value = 756.1207216813057
$ws.Range("R12").Value = 6805.08649513175
$ws.Range("S12").Value = 0.102252504700543
$ws.Range("T12").Value = 0.102252504700543
$ws.Range("G13").Value = 12.37490333333333
$ws.Range("H13").Value = 37.12471
$ws.Range("I13").Value = 0.2057223442826371
$ws.Range("J13").Value = 0.2057223442826371
$ws.Range("M13").Value = 0.1059033333333333
$ws.Range("N13").Value = 0.31771
$ws.Range("O13").Value = 0.0008614950672987739
$ws.Range("P13").Value = 0.0008614950672987739
$ws.Range("Q13").Value = 1.310543512677778
$ws.Range("R13").Value = 11.7948916141
$ws.Range("S13").Value = 0.000177228784832632
$ws.Range("T13").Value = 0.000177228784832632
